$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 204 (this shifts the existing rows 204..318 down to 205..319,
# preserving their values and formats untouched).
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new weekly record.
# (Same dimension/category columns as the former row 204, with updated
# date / volume / price figures.)
$ws.Range("A204").Value = 11
$ws.Range("B204").Value = "Vega Monumental Concepción"
$ws.Range("C204").Value = "Bíobío"
$ws.Range("D204").Value = 45086
$ws.Range("E204").Value = 8
$ws.Range("F204").Value = 100112040
$ws.Range("G204").Value = "Cilantro"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 50
$ws.Range("K204").Value = 7500
$ws.Range("L204").Value = 8000
$ws.Range("M204").Value = 7800
$ws.Range("N204").Value = "$/caja 36 atados"
$ws.Range("O204").Value = "Región Metropolitana"
$ws.Range("P204").Value = 217
$ws.Range("Q204").Value = 36
$ws.Range("R204").Value = "Hortaliza"
